$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 18.82499266666667
$ws.Cells.Item(2, 8).Value = 56.474978
$ws.Cells.Item(2, 9).Value = 0.06886869772378311
$ws.Cells.Item(2, 10).Value = 0.0688686977237831
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 46.33695966666667
$ws.Cells.Item(2, 14).Value = 139.010879
$ws.Cells.Item(2, 15).Value = 0.1993490803952133
$ws.Cells.Item(2, 16).Value = 0.1993490803952133
$ws.Cells.Item(2, 17).Value = 872.292925920629
$ws.Cells.Item(2, 18).Value = 7850.636333285662
$ws.Cells.Item(2, 19).Value = 0.01372891155925208
$ws.Cells.Item(2, 20).Value = 0.01372891155925208

# Row 3
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 18.82499266666667
$ws.Cells.Item(3, 8).Value = 56.474978
$ws.Cells.Item(3, 9).Value = 0.06886869772378311
$ws.Cells.Item(3, 10).Value = 0.0688686977237831
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 84.50960033333332
$ws.Cells.Item(3, 14).Value = 253.528801
$ws.Cells.Item(3, 15).Value = 0.3635739425333109
$ws.Cells.Item(3, 16).Value = 0.3635739425333109
$ws.Cells.Item(3, 17).Value = 1590.892606537931
$ws.Cells.Item(3, 18).Value = 14318.03345884138
$ws.Cells.Item(3, 19).Value = 0.02503886394857068
$ws.Cells.Item(3, 20).Value = 0.02503886394857067

# Row 4
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 18.82499266666667
$ws.Cells.Item(4, 8).Value = 56.474978
$ws.Cells.Item(4, 9).Value = 0.06886869772378311
$ws.Cells.Item(4, 10).Value = 0.0688686977237831
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 72.52790466666666
$ws.Cells.Item(4, 14).Value = 217.583714
$ws.Cells.Item(4, 15).Value = 0.3120267536390091
$ws.Cells.Item(4, 16).Value = 0.3120267536390091
$ws.Cells.Item(4, 17).Value = 1365.337273478699
$ws.Cells.Item(4, 18).Value = 12288.03546130829
$ws.Cells.Item(4, 19).Value = 0.02148887617809826
$ws.Cells.Item(4, 20).Value = 0.02148887617809825

# Row 5
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 18.82499266666667
$ws.Cells.Item(5, 8).Value = 56.474978
$ws.Cells.Item(5, 9).Value = 0.06886869772378311
$ws.Cells.Item(5, 10).Value = 0.0688686977237831
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 29.06683666666666
$ws.Cells.Item(5, 14).Value = 87.20050999999998
$ws.Cells.Item(5, 15).Value = 0.1250502234324667
$ws.Cells.Item(5, 16).Value = 0.1250502234324667
$ws.Cells.Item(5, 17).Value = 547.1829870931977
$ws.Cells.Item(5, 18).Value = 4924.646883838779
$ws.Cells.Item(5, 19).Value = 0.008612046037862092
$ws.Cells.Item(5, 20).Value = 0.00861204603786209

# Row 6
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 121.8208923333333
$ws.Cells.Item(6, 8).Value = 365.462677
$ws.Cells.Item(6, 9).Value = 0.4456653109566078
$ws.Cells.Item(6, 10).Value = 0.4456653109566078
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 46.33695966666667
$ws.Cells.Item(6, 14).Value = 139.010879
$ws.Cells.Item(6, 15).Value = 0.1993490803952133
$ws.Cells.Item(6, 16).Value = 0.1993490803952133
$ws.Cells.Item(6, 17).Value = 5644.80977460701
$ws.Cells.Item(6, 18).Value = 50803.28797146308
$ws.Cells.Item(6, 19).Value = 0.08884296990324655
$ws.Cells.Item(6, 20).Value = 0.08884296990324653

# Row 7
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 121.8208923333333
$ws.Cells.Item(7, 8).Value = 365.462677
$ws.Cells.Item(7, 9).Value = 0.4456653109566078
$ws.Cells.Item(7, 10).Value = 0.4456653109566078
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 84.50960033333332
$ws.Cells.Item(7, 14).Value = 253.528801
$ws.Cells.Item(7, 15).Value = 0.3635739425333109
$ws.Cells.Item(7, 16).Value = 0.3635739425333109
$ws.Cells.Item(7, 17).Value = 10295.03492334003
$ws.Cells.Item(7, 18).Value = 92655.31431006029
$ws.Cells.Item(7, 19).Value = 0.1620322941548279
$ws.Cells.Item(7, 20).Value = 0.1620322941548278

# Row 8
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 121.8208923333333
$ws.Cells.Item(8, 8).Value = 365.462677
$ws.Cells.Item(8, 9).Value = 0.4456653109566078
$ws.Cells.Item(8, 10).Value = 0.4456653109566078
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 72.52790466666666
$ws.Cells.Item(8, 14).Value = 217.583714
$ws.Cells.Item(8, 15).Value = 0.3120267536390091
$ws.Cells.Item(8, 16).Value = 0.3120267536390091
$ws.Cells.Item(8, 17).Value = 8835.414065560264
$ws.Cells.Item(8, 18).Value = 79518.72659004238
$ws.Cells.Item(8, 19).Value = 0.1390595001873098
$ws.Cells.Item(8, 20).Value = 0.1390595001873098

# Row 9
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 121.8208923333333
$ws.Cells.Item(9, 8).Value = 365.462677
$ws.Cells.Item(9, 9).Value = 0.4456653109566078
$ws.Cells.Item(9, 10).Value = 0.4456653109566078
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 29.06683666666666
$ws.Cells.Item(9, 14).Value = 87.20050999999998
$ws.Cells.Item(9, 15).Value = 0.1250502234324667
$ws.Cells.Item(9, 16).Value = 0.1250502234324667
$ws.Cells.Item(9, 17).Value = 3540.947980040585
$ws.Cells.Item(9, 18).Value = 31868.53182036527
$ws.Cells.Item(9, 19).Value = 0.05573054671122358
$ws.Cells.Item(9, 20).Value = 0.05573054671122357

# Row 10
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 87.673585
$ws.Cells.Item(10, 8).Value = 263.020755
$ws.Cells.Item(10, 9).Value = 0.3207419907481189
$ws.Cells.Item(10, 10).Value = 0.3207419907481188
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 46.33695966666667
$ws.Cells.Item(10, 14).Value = 139.010879
$ws.Cells.Item(10, 15).Value = 0.1993490803952133
$ws.Cells.Item(10, 16).Value = 0.1993490803952133
$ws.Cells.Item(10, 17).Value = 4062.527371977072
$ws.Cells.Item(10, 18).Value = 36562.74634779364
$ws.Cells.Item(10, 19).Value = 0.06393962089976751
$ws.Cells.Item(10, 20).Value = 0.0639396208997675

# Row 11
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 87.673585
$ws.Cells.Item(11, 8).Value = 263.020755
$ws.Cells.Item(11, 9).Value = 0.3207419907481189
$ws.Cells.Item(11, 10).Value = 0.3207419907481188
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 13).Value = 84.50960033333332
$ws.Cells.Item(11, 14).Value = 253.528801
$ws.Cells.Item(11, 15).Value = 0.3635739425333109
$ws.Cells.Item(11, 16).Value = 0.3635739425333109
$ws.Cells.Item(11, 17).Value = 7409.259628140528
$ws.Cells.Item(11, 18).Value = 66683.33665326475
$ws.Cells.Item(11, 19).Value = 0.1166134301122763
$ws.Cells.Item(11, 20).Value = 0.1166134301122763

# Row 12
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 7).Value = 87.673585
$ws.Cells.Item(12, 8).Value = 263.020755
$ws.Cells.Item(12, 9).Value = 0.3207419907481189
$ws.Cells.Item(12, 10).Value = 0.3207419907481188
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 72.52790466666666
$ws.Cells.Item(12, 14).Value = 217.583714
$ws.Cells.Item(12, 15).Value = 0.3120267536390091
$ws.Cells.Item(12, 16).Value = 0.3120267536390091
$ws.Cells.Item(12, 17).Value = 6358.781414664896
$ws.Cells.Item(12, 18).Value = 57229.03273198407
$ws.Cells.Item(12, 19).Value = 0.1000800821288486
$ws.Cells.Item(12, 20).Value = 0.1000800821288486

# Row 13
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 7).Value = 87.673585
$ws.Cells.Item(13, 8).Value = 263.020755
$ws.Cells.Item(13, 9).Value = 0.3207419907481189
$ws.Cells.Item(13, 10).Value = 0.3207419907481188
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 13).Value = 29.06683666666666
$ws.Cells.Item(13, 14).Value = 87.20050999999998
$ws.Cells.Item(13, 15).Value = 0.1250502234324667
$ws.Cells.Item(13, 16).Value = 0.1250502234324667
$ws.Cells.Item(13, 17).Value = 2548.393775176116
$ws.Cells.Item(13, 18).Value = 22935.54397658505
$ws.Cells.Item(13, 19).Value = 0.04010885760722645
$ws.Cells.Item(13, 20).Value = 0.04010885760722644

# Row 14
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 7).Value = 45.02666966666666
$ws.Cells.Item(14, 8).Value = 135.080009
$ws.Cells.Item(14, 9).Value = 0.1647240005714903
$ws.Cells.Item(14, 10).Value = 0.1647240005714903
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 13).Value = 46.33695966666667
$ws.Cells.Item(14, 14).Value = 139.010879
$ws.Cells.Item(14, 15).Value = 0.1993490803952133
$ws.Cells.Item(14, 16).Value = 0.1993490803952133
$ws.Cells.Item(14, 17).Value = 2086.398976268656
$ws.Cells.Item(14, 18).Value = 18777.59078641791
$ws.Cells.Item(14, 19).Value = 0.03283757803294718
$ws.Cells.Item(14, 20).Value = 0.03283757803294717

# Row 15
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 7).Value = 45.02666966666666
$ws.Cells.Item(15, 8).Value = 135.080009
$ws.Cells.Item(15, 9).Value = 0.1647240005714903
$ws.Cells.Item(15, 10).Value = 0.1647240005714903
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 13).Value = 84.50960033333332
$ws.Cells.Item(15, 14).Value = 253.528801
$ws.Cells.Item(15, 15).Value = 0.3635739425333109
$ws.Cells.Item(15, 16).Value = 0.3635739425333109
$ws.Cells.Item(15, 17).Value = 3805.185857871023
$ws.Cells.Item(15, 18).Value = 34246.6727208392
$ws.Cells.Item(15, 19).Value = 0.05988935431763608
$ws.Cells.Item(15, 20).Value = 0.05988935431763607

# Row 16
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 7).Value = 45.02666966666666
$ws.Cells.Item(16, 8).Value = 135.080009
$ws.Cells.Item(16, 9).Value = 0.1647240005714903
$ws.Cells.Item(16, 10).Value = 0.1647240005714903
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 13).Value = 72.52790466666666
$ws.Cells.Item(16, 14).Value = 217.583714
$ws.Cells.Item(16, 15).Value = 0.3120267536390091
$ws.Cells.Item(16, 16).Value = 0.3120267536390091
$ws.Cells.Item(16, 17).Value = 3265.690005041491
$ws.Cells.Item(16, 18).Value = 29391.21004537342
$ws.Cells.Item(16, 19).Value = 0.05139829514475239
$ws.Cells.Item(16, 20).Value = 0.05139829514475237

# Row 17
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 7).Value = 45.02666966666666
$ws.Cells.Item(17, 8).Value = 135.080009
$ws.Cells.Item(17, 9).Value = 0.1647240005714903
$ws.Cells.Item(17, 10).Value = 0.1647240005714903
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 13).Value = 29.06683666666666
$ws.Cells.Item(17, 14).Value = 87.20050999999998
$ws.Cells.Item(17, 15).Value = 0.1250502234324667
$ws.Cells.Item(17, 16).Value = 0.1250502234324667
$ws.Cells.Item(17, 17).Value = 1308.782852844954
$ws.Cells.Item(17, 18).Value = 11779.04567560459
$ws.Cells.Item(17, 19).Value = 0.02059877307615464
$ws.Cells.Item(17, 20).Value = 0.02059877307615464
